$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new column I ("Valid")
$ws.Range("I1").Value = "Valid"

# Existing rows 2-13: add Valid = 1
$ws.Range("I2").Value = 1
$ws.Range("I3").Value = 1
$ws.Range("I4").Value = 1
$ws.Range("I5").Value = 1
$ws.Range("I6").Value = 1
$ws.Range("I7").Value = 1
$ws.Range("I8").Value = 1
$ws.Range("I9").Value = 1
$ws.Range("I10").Value = 1
$ws.Range("I11").Value = 1
$ws.Range("I12").Value = 1
$ws.Range("I13").Value = 1

# Rows 14-22: full new data rows
$ws.Range("B14").Value = "G"
$ws.Range("C14").Value = 3
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 9
$ws.Range("G14").Value = 7
$ws.Range("H14").Value = 1
$ws.Range("I14").Value = 1

$ws.Range("B15").Value = "F"
$ws.Range("C15").Value = 14
$ws.Range("D15").Value = 8
$ws.Range("E15").Value = 10
$ws.Range("F15").Value = 19
$ws.Range("G15").Value = 18
$ws.Range("H15").Value = 15
$ws.Range("I15").Value = 0

$ws.Range("B16").Value = "G"
$ws.Range("C16").Value = 14
$ws.Range("D16").Value = 10
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 4
$ws.Range("G16").Value = 11
$ws.Range("H16").Value = 13
$ws.Range("I16").Value = 0

$ws.Range("B17").Value = "F"
$ws.Range("C17").Value = 14
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = 10
$ws.Range("F17").Value = 9
$ws.Range("G17").Value = 13
$ws.Range("H17").Value = 17
$ws.Range("I17").Value = 1

$ws.Range("B18").Value = "G"
$ws.Range("C18").Value = 16
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 1
$ws.Range("F18").Value = 14
$ws.Range("G18").Value = 15
$ws.Range("H18").Value = 7
$ws.Range("I18").Value = 1

$ws.Range("B19").Value = "F"
$ws.Range("C19").Value = 13
$ws.Range("D19").Value = 2
$ws.Range("E19").Value = 10
$ws.Range("F19").Value = 2
$ws.Range("G19").Value = 14
$ws.Range("H19").Value = 6
$ws.Range("I19").Value = 1

$ws.Range("B20").Value = "G"
$ws.Range("C20").Value = 8
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = 12
$ws.Range("F20").Value = 13
$ws.Range("G20").Value = 11
$ws.Range("H20").Value = 2
$ws.Range("I20").Value = 1

$ws.Range("B21").Value = "F"
$ws.Range("C21").Value = 13
$ws.Range("D21").Value = 4
$ws.Range("E21").Value = 5
$ws.Range("F21").Value = 10
$ws.Range("G21").Value = 10
$ws.Range("H21").Value = 5
$ws.Range("I21").Value = 1

$ws.Range("B22").Value = "G"
$ws.Range("C22").Value = 8
$ws.Range("D22").Value = 2
$ws.Range("E22").Value = 11
$ws.Range("F22").Value = 7
$ws.Range("G22").Value = 11
$ws.Range("H22").Value = 13
$ws.Range("I22").Value = 1

# Rows 23-32: only Condition (B) and Valid (I)
$ws.Range("B23").Value = "F"
$ws.Range("I23").Value = 1
$ws.Range("B24").Value = "G"
$ws.Range("I24").Value = 1
$ws.Range("B25").Value = "F"
$ws.Range("I25").Value = 1
$ws.Range("B26").Value = "G"
$ws.Range("I26").Value = 1
$ws.Range("B27").Value = "F"
$ws.Range("I27").Value = 1
$ws.Range("B28").Value = "G"
$ws.Range("I28").Value = 1
$ws.Range("B29").Value = "F"
$ws.Range("I29").Value = 1
$ws.Range("B30").Value = "G"
$ws.Range("I30").Value = 1
$ws.Range("B31").Value = "F"
$ws.Range("I31").Value = 1
$ws.Range("B32").Value = "G"
$ws.Range("I32").Value = 1

# Update selection to match target
$ws.Range("O16").Select()
